{"js": "// The document contains four \"use case\" tables, each with a row labeled\n// \"Post-Conditions\" and a row labeled \"Preconditions\". The edit swaps the\n// labels (the commit message says \"swapped preconditions and\n// post-conditions\"): every \"Post-Conditions\" cell becomes \"Preconditions\"\n// and every \"Preconditions\" cell becomes \"Post-conditions\" (note the\n// lower-case \"c\" used by the author for the post-condition label).\n\n// Search for the two exact labels first (before any mutation), so the\n// search results/ranges are computed against the original, unmodified\n// text and are not confused by text we are about to insert.\nconst postConditionRanges = context.document.body.search(\"Post-Conditions\", {\n  matchCase: true,\n  matchWholeWord: true,\n});\npostConditionRanges.load(\"items\");\n\nconst preConditionRanges = context.document.body.search(\"Preconditions\", {\n  matchCase: true,\n  matchWholeWord: true,\n});\npreConditionRanges.load(\"items\");\n\nawait context.sync();\n\nconst postCount = postConditionRanges.items.length;\nconst preCount = preConditionRanges.items.length;\n\n// \"Post-Conditions\" -> \"Preconditions\"\nfor (let i = 0; i < postCount; i++) {\n  postConditionRanges.items[i].insertText(\"Preconditions\", Word.InsertLocation.replace);\n}\n\n// \"Preconditions\" -> \"Post-conditions\"\nfor (let i = 0; i < preCount; i++) {\n  preConditionRanges.items[i].insertText(\"Post-conditions\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document contains four \"use case\" tables, each with a row labeled\n# \"Post-Conditions\" and a row labeled \"Preconditions\". The edit swaps the\n# labels (commit message: \"swapped preconditions and post-conditions\"):\n# every \"Post-Conditions\" cell becomes \"Preconditions\" and every\n# \"Preconditions\" cell becomes \"Post-conditions\" (note the lower-case \"c\"\n# used by the author for the post-condition label).\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText {\n    param(\n        $Document,\n        [string]$SearchText,\n        [string]$ReplaceText\n    )\n\n    $find = $Document.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $SearchText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.Replacement.Text = $ReplaceText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($SearchText, $true, $true, $false, $false, $false, $true, 1, $false, $ReplaceText, 2) | Out-Null\n}\n\n# A plain text-for-text swap would clobber itself (the second replace\n# would also match the text produced by the first replace), so route the\n# swap through a placeholder that cannot already occur in the document.\n$placeholder = \"###SWAP_POST_CONDITIONS###\"\n\nReplace-AllText $d \"Post-Conditions\" $placeholder\nReplace-AllText $d \"Preconditions\" \"Post-conditions\"\nReplace-AllText $d $placeholder \"Preconditions\"\n"}
